# Auto-generated edit script: apply numeric updates to Golem_Profits sheets
# Source: scheduled runner market-data refresh (chore: update Sheets via scheduled runner)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 2
$ws.Range("H2").Value = 729
$ws.Range("J2").Value = 729
$ws.Range("L2").Value = 729
$ws.Range("N2").Value = -955
# row 26
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
# row 32
$ws.Range("H32").Value = 2133
$ws.Range("I32").Value = 1199.5
$ws.Range("J32").Value = 4000
$ws.Range("K32").Value = 1199.5
$ws.Range("L32").Value = 4000
$ws.Range("M32").Value = -873.5
$ws.Range("N32").Value = -4652
# row 39
$ws.Range("H39").Value = 52.842106
$ws.Range("I39").Value = 55.941177
$ws.Range("J39").Value = 26.5
$ws.Range("K39").Value = 167.823531
$ws.Range("L39").Value = 79.5
$ws.Range("M39").Value = 128.176469
$ws.Range("N39").Value = -671.5
# row 48
$ws.Range("H48").Value = 2000
$ws.Range("I48").Value = 2000
$ws.Range("K48").Value = 6000
$ws.Range("M48").Value = -5708
# row 51
$ws.Range("H51").Value = 6999
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 6999
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 6999
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -7967
# row 54
$ws.Range("H54").Value = 38997
$ws.Range("J54").Value = 38997
$ws.Range("L54").Value = 38997
$ws.Range("N54").Value = -39969
# row 56
$ws.Range("H56").Value = 2000
$ws.Range("I56").Value = 2000
$ws.Range("K56").Value = 6000
$ws.Range("M56").Value = -5466
# row 132
$ws.Range("H132").Value = 2489.6
$ws.Range("I132").Value = 2455.3333
$ws.Range("K132").Value = 7365.999899999999
$ws.Range("M132").Value = -4835.999899999999
# row 137
$ws.Range("H137").Value = 2886
$ws.Range("I137").Value = 3475.25
$ws.Range("J137").Value = 2100.3333
$ws.Range("K137").Value = 10425.75
$ws.Range("L137").Value = 6300.999899999999
$ws.Range("M137").Value = -7875.75
$ws.Range("N137").Value = -11400.9999
# row 138
$ws.Range("H138").Value = 3543
$ws.Range("I138").Value = 3543
$ws.Range("K138").Value = 10629
$ws.Range("M138").Value = -5489

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 1728.5264
$ws.Range("I32").Value = 1663.9333
$ws.Range("K32").Value = 1663.9333
$ws.Range("M32").Value = -1376.9333
# row 74
$ws.Range("H74").Value = 447589.88
$ws.Range("I74").Value = 574387.1
$ws.Range("K74").Value = 574387.1
$ws.Range("M74").Value = -573513.1
# row 77
$ws.Range("H77").Value = 447589.88
$ws.Range("I77").Value = 574387.1
$ws.Range("K77").Value = 2871935.5
$ws.Range("M77").Value = -2867567.5
# row 88
$ws.Range("H88").Value = 1625.3077
$ws.Range("J88").Value = 966.2727
$ws.Range("L88").Value = 966.2727
$ws.Range("N88").Value = -1778.2727
# row 91
$ws.Range("H91").Value = 1625.3077
$ws.Range("J91").Value = 966.2727
$ws.Range("L91").Value = 966.2727
$ws.Range("N91").Value = -3774.2727
# row 122
$ws.Range("H122").Value = 11099.286
$ws.Range("I122").Value = 13333.333
$ws.Range("J122").Value = 9423.75
$ws.Range("K122").Value = 39999.999
$ws.Range("L122").Value = 28271.25
$ws.Range("M122").Value = -37549.999
$ws.Range("N122").Value = -33171.25
# row 132
$ws.Range("H132").Value = 4420
$ws.Range("I132").Value = 4420
$ws.Range("K132").Value = 13260
$ws.Range("M132").Value = -10730

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 86
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
# row 89
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
# row 100
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
# row 134
$ws.Range("H134").Value = 4216.25
$ws.Range("I134").Value = 4216.25
$ws.Range("K134").Value = 12648.75
$ws.Range("M134").Value = -10113.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 1428.4
$ws.Range("I31").Value = 1459.7142
$ws.Range("J31").Value = 990
$ws.Range("K31").Value = 1459.7142
$ws.Range("L31").Value = 990
$ws.Range("M31").Value = -1164.7142
$ws.Range("N31").Value = -1580
# row 34
$ws.Range("H34").Value = 1428.4
$ws.Range("I34").Value = 1459.7142
$ws.Range("J34").Value = 990
$ws.Range("K34").Value = 1459.7142
$ws.Range("L34").Value = 990
$ws.Range("M34").Value = -1257.7142
$ws.Range("N34").Value = -1394
# row 81
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
# row 84
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
# row 87
$ws.Range("H87").Value = 45000
$ws.Range("J87").Value = 45000
$ws.Range("L87").Value = 45000
$ws.Range("N87").Value = -47372
# row 90
$ws.Range("H90").Value = 45000
$ws.Range("J90").Value = 45000
$ws.Range("L90").Value = 135000
$ws.Range("N90").Value = -146856

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 46
$ws.Range("H46").Value = 3963.8333
$ws.Range("J46").Value = 5495
$ws.Range("L46").Value = 16485
$ws.Range("N46").Value = -16667
# row 50
$ws.Range("H50").Value = 458.33334
$ws.Range("I50").Value = 458.33334
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 1375.00002
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -894.0000199999999
$ws.Range("N50").ClearContents()
# row 53
$ws.Range("H53").Value = 458.33334
$ws.Range("I53").Value = 458.33334
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 1375.00002
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -894.0000199999999
$ws.Range("N53").ClearContents()
# row 54
$ws.Range("H54").Value = 4900
$ws.Range("J54").Value = 4900
$ws.Range("L54").Value = 14700
$ws.Range("N54").Value = -15818
# row 87
$ws.Range("H87").Value = 3506.5
$ws.Range("I87").Value = 3506.5
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 10519.5
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -9271.5
$ws.Range("N87").ClearContents()
# row 90
$ws.Range("H90").Value = 3506.5
$ws.Range("I90").Value = 3506.5
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 31558.5
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -25318.5
$ws.Range("N90").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 122
$ws.Range("H122").Value = 5214.857
$ws.Range("I122").Value = 4401.2
$ws.Range("J122").Value = 7249
$ws.Range("K122").Value = 13203.6
$ws.Range("L122").Value = 21747
$ws.Range("M122").Value = -10753.6
$ws.Range("N122").Value = -26647
# row 132
$ws.Range("H132").Value = 2713.0715
$ws.Range("J132").Value = 3797.6
$ws.Range("L132").Value = 11392.8
$ws.Range("N132").Value = -16452.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 93
$ws.Range("H93").Value = 1370.2858
$ws.Range("I93").Value = 1298.8334
$ws.Range("K93").Value = 1298.8334
$ws.Range("M93").Value = -50.83339999999998
# row 136
$ws.Range("H136").Value = 1274750.8
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
# row 139
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
# row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
# row 132
$ws.Range("H132").Value = 2307.7568
$ws.Range("I132").Value = 2296.6287
$ws.Range("K132").Value = 6889.886100000001
$ws.Range("M132").Value = -4359.886100000001
